$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values per data row (row index -> column index -> value)
$data = @{
  2 = @{ 2=20.75921928337546; 3=5.026124265168812; 4=13.55037017520021; 5=13.62159553492157; 7=67.81168018883096; 8=24.08555380105708; 10=8.311535738884389; 11=16.36618823416302; 12=12.57644374224787; 13=19.53977315126451; 14=26.16797549182537 }
  3 = @{ 2=20.66871965249276; 3=4.915149957352502; 4=13.55478957553749; 5=13.64418459480357; 7=67.79503116842854; 8=24.11951767727522; 10=8.313027636936582; 11=16.30584594711877; 12=12.593853693817; 13=19.54719799014584; 14=26.21839251963568 }
  4 = @{ 2=20.61766586050066; 3=4.844700097549159; 4=13.55967803221928; 5=13.65933404124704; 7=67.79715272644356; 8=24.14356877746898; 10=8.31401435598746; 11=16.27229001793161; 12=12.60607095595971; 13=19.55483893966151; 14=26.25126866315533 }
  5 = @{ 2=20.59801264618386; 3=4.815428077898314; 4=13.56221769196501; 5=13.66582984658823; 7=67.80112267993186; 8=24.15417354632561; 10=8.314434264180399; 11=16.25950502840369; 12=12.61143407375395; 13=19.55872889218459; 14=26.26514934478741 }
  6 = @{ 2=20.5948192582959; 3=4.810534019304832; 4=13.56267248882314; 5=13.6669279497702; 7=67.80196945718319; 8=24.1559829972378; 10=8.31450506652571; 11=16.25743609301742; 12=12.61234784798184; 13=19.55942172688246; 14=26.26748343502042 }
  7 = @{ 2=20.61739612659347; 3=4.844307580258447; 4=13.55971006509129; 5=13.65942034037708; 7=67.79719369238479; 8=24.14370854310231; 10=8.314019946844963; 11=16.27211398103526; 12=12.60614172759288; 13=19.55488825671304; 14=26.2514539048811 }
  8 = @{ 2=20.72708873767442; 3=4.988350468543258; 4=13.55144302928529; 5=13.62911894755931; 7=67.80337882286192; 8=24.09660088130286; 10=8.312035496393765; 11=16.34466361123498; 12=12.58212990493362; 13=19.5416943322914; 14=26.18496118101341 }
  9 = @{ 2=20.97720517728469; 3=5.25174870322224; 4=13.55245337240626; 5=13.57982963860028; 7=67.91334001766096; 8=24.02959833805723; 10=8.30870323134106; 11=16.51414767920727; 12=12.5471452299452; 13=19.54021612691539; 14=26.06977453488471 }
  10 = @{ 2=21.18114133985088; 3=5.43274867819324; 4=13.56363916554069; 5=13.54976359821744; 7=68.0535402252129; 8=23.9958494387503; 10=8.306593733030256; 11=16.65450688900323; 12=12.52879578877015; 13=19.55391066447495; 14=25.9943773676688 }
  11 = @{ 2=21.27802273260266; 3=5.512206127245975; 4=13.57098000991989; 5=13.53741413319904; 7=68.13014096485469; 8=23.983857149044; 10=8.305707155959539; 11=16.72161991370907; 12=12.52203917462886; 13=19.56332569160606; 14=25.9620734677167 }
  12 = @{ 2=21.31527522242205; 3=5.541867951446088; 4=13.57408200088203; 5=13.53292811656857; 7=68.16098179487715; 8=23.97979897608284; 10=8.305381900234748; 11=16.74748662957059; 12=12.51970879289895; 13=19.56734640743417; 14=25.95012709297317 }
  13 = @{ 2=21.30722749104892; 3=5.535498929467793; 4=13.57339963234627; 5=13.53388579736597; 7=68.15425829052684; 8=23.98065149617788; 10=8.305451484643445; 11=16.7418959201759; 12=12.52020054186219; 13=19.56646026222936; 14=25.95268722677095 }
  14 = @{ 2=21.28107633055743; 3=5.514655050267341; 4=13.57122876629695; 5=13.53704125155364; 7=68.13264157039769; 8=23.98351360099647; 10=8.305680187286779; 11=16.72373901433767; 12=12.52184288240866; 13=19.56364736623418; 14=25.96108489450504 }
  15 = @{ 2=21.26513087864866; 3=5.501831604855577; 4=13.56994095023167; 5=13.53899884759688; 7=68.11963919052937; 8=23.98532962340921; 10=8.305821637053281; 11=16.71267580067325; 12=12.52287856553596; 13=19.56198361765482; 14=25.96626599554216 }
  16 = @{ 2=21.17489046656186; 3=5.427496868187757; 4=13.56320460542467; 5=13.55059734220586; 7=68.04879138839233; 8=23.99670076930055; 10=8.306653140027109; 11=16.65018514029525; 12=12.52926931268921; 13=19.55335925117942; 14=25.99652860492094 }
  17 = @{ 2=21.12056578594715; 3=5.381147735325347; 4=13.55964780455161; 5=13.55805237471468; 7=68.00860613182373; 8=24.00453716532804; 10=8.307181925720073; 11=16.61267319565287; 12=12.53359686519898; 13=19.54888255456486; 14=26.01560431127338 }
  18 = @{ 2=21.08970808778128; 3=5.354218522388245; 4=13.55781414583511; 5=13.56246531256125; 7=67.98670030270098; 8=24.00936075180843; 10=8.307492946422546; 11=16.59140549084307; 12=12.5362357116557; 13=19.54660766581393; 14=26.02676392267142 }
  19 = @{ 2=21.079327669755; 3=5.345054716967768; 4=13.5572297811799; 5=13.56398094058197; 7=67.97949106565488; 8=24.01104826290351; 10=8.30759943482529; 11=16.5842580377824; 12=12.5371549154009; 13=19.54588902567471; 14=26.0305746379772 }
  20 = @{ 2=21.12630872141963; 3=5.386109763025141; 4=13.56000449035376; 5=13.55724584008611; 7=68.01275900006482; 8=24.00367023259061; 10=8.307124924053602; 11=16.61663463016237; 12=12.53312069471557; 13=19.54932807680067; 14=26.01355423935192 }
  21 = @{ 2=21.28874241694648; 3=5.520789087755329; 4=13.57185767326462; 5=13.53610925295734; 7=68.13894124531176; 8=23.98265982343332; 10=8.305612727827857; 11=16.72905999489; 12=12.52135429828819; 13=19.56446124319553; 14=25.95861052605591 }
  22 = @{ 2=21.39818659172687; 3=5.606315826341499; 4=13.58148134590024; 5=13.52340518387565; 7=68.23209153758238; 8=23.97174372591945; 10=8.304685443173275; 11=16.80516553489798; 12=12.51499413549603; 13=19.57700496002814; 14=25.92437087428498 }
  23 = @{ 2=21.33948228750187; 3=5.560900744598149; 4=13.57617388919993; 5=13.53008418359257; 7=68.1814017545579; 8=23.97731232017041; 10=8.305174779264389; 11=16.76431173492301; 12=12.51826717783812; 13=19.57006828892349; 14=25.94249261877776 }
  24 = @{ 2=21.12371117248705; 3=5.383867309232379; 4=13.55984257473886; 5=13.5576100789524; 7=68.01087775819754; 8=24.0040611811973; 10=8.307150672656055; 11=16.61484273480974; 12=12.53333550137509; 13=19.54912572517748; 14=26.01448047641751 }
  25 = @{ 2=20.90590922956817; 3=5.182630774208008; 4=13.550341288954; 5=13.59208209304224; 7=67.87313638124971; 8=24.04500644756627; 10=8.309545057141449; 11=16.46545802366552; 12=12.55531621353394; 13=19.53801278059839; 14=26.09931184085578 }
}

foreach ($r in $data.Keys) {
  foreach ($c in $data[$r].Keys) {
    $ws.Cells.Item($r, $c).Value = $data[$r][$c]
  }
}

Write-Output "Updated $($data.Keys.Count) rows"
